$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.267.24'
$ws.Range("E2").Value = '  -1.21%  '
$ws.Range("D3").Value = '3.338.29'
$ws.Range("E3").Value = '  +2.45%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.03'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.32%  '
$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.603'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.52%  '
$ws.Range("E9").Value = '  -0.16%  '
$ws.Range("E10").Value = '  +0.82%  '
$ws.Range("E11").Value = '  +0.03%  '
$ws.Range("D12").Value = '3.923.35'
$ws.Range("E12").Value = '  +2.51%  '
$ws.Range("E13").Value = '  -0.76%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.30'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.33%  '
$ws.Range("D15").Value = '67.461.16'
$ws.Range("E15").Value = '  -0.91%  '
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").Value = '3.344.32'
$ws.Range("E17").Value = '  +2.55%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '445.10'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +7.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.58'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.65'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.82%  '
$ws.Range("E21").Value = '  +2.87%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.98'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.05%  '
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("D24").Value = '3.491.47'
$ws.Range("E24").Value = '  +2.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.511'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.95%  '
$ws.Range("E26").Value = '  +3.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.193'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.03'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.56%  '
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("E30").Value = '  +1.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.92'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.85%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.32'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.92%  '
$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.23'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.77'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.75%  '
$ws.Range("E36").Value = '  +4.38%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '161.72'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.57%  '
$ws.Range("E38").Value = '  +4.04%  '
$ws.Range("E39").Value = '  -1.50%  '
$ws.Range("D40").Value = '2.833.71'
$ws.Range("E40").Value = '  +8.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.791'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.46'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.24'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.31'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.65%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0671'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '24.50'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.97%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.35'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.99%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '322.93'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.45%  '
$ws.Range("E49").Value = '  +0.77%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.986'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.95%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '30.93'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.76%  '
